$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "main"

# Clear the Slope label and formula cells
$ws.Range("K5").ClearContents()
$ws.Range("K6").ClearContents()

# Update selection
$ws.Range("P7").Select()
